$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These D-column cells get a new value that *looks* purely numeric
# ("228.71", "1.00", "0.0920", ...). The source data is plain text
# (European-grouped price strings elsewhere use two dots, e.g.
# "38.812.04", which naturally stay text) and must stay text here too,
# so trailing zeros / exact formatting survive. Force the cell format to
# Text before assigning so the engine doesn't reinterpret the string as
# a number.
$textForceCells = @(
    "D5","D7","D10","D11","D13","D14","D15","D19","D22","D24","D26",
    "D28","D29","D30","D36","D37","D38","D39","D40","D43","D45","D46","D50"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    2  = @{ D = "38.812.04"; E = "  +2.83%  " }
    3  = @{ D = "2.090.48";  E = "  +2.13%  " }
    4  = @{ E = "  -0.03%  " }
    5  = @{ D = "228.71";    E = "  +0.49%  " }
    6  = @{ E = "  +0.82%  " }
    7  = @{ D = "60.21";     E = "  +0.54%  " }
    8  = @{ E = "  -0.03%  " }
    9  = @{ E = "  +2.00%  " }
    10 = @{ D = "0.0839";    E = "  +0.20%  " }
    11 = @{ D = "0.103";     E = "  -0.90%  " }
    12 = @{ D = "2.403.72";  E = "  +2.30%  " }
    13 = @{ D = "14.97";     E = "  +4.13%  " }
    14 = @{ D = "21.92";     E = "  +2.58%  " }
    15 = @{ D = "0.797";     E = "  +4.14%  " }
    16 = @{ E = "  -0.80%  " }
    17 = @{ D = "2.095.87";  E = "  +1.65%  " }
    18 = @{ D = "38.729.31"; E = "  +2.57%  " }
    19 = @{ D = "71.56";     E = "  +3.11%  " }
    20 = @{ E = "  +2.42%  " }
    21 = @{ E = "  +1.07%  " }
    22 = @{ D = "226.77";    E = "  +1.93%  " }
    24 = @{ D = "2.38";      E = "  -0.17%  " }
    25 = @{ E = "  +2.72%  " }
    26 = @{ D = "171.10";    E = "  +1.34%  " }
    27 = @{ E = "  +2.09%  " }
    28 = @{ D = "0.141";     E = "  +10.34%  " }
    29 = @{ D = "1.47";      E = "  +13.83%  " }
    30 = @{ D = "19.17";     E = "  +2.19%  " }
    31 = @{ E = "  +0.98%  " }
    32 = @{ E = "  +5.14%  " }
    33 = @{ E = "  +2.80%  " }
    34 = @{ E = "  +4.02%  " }
    35 = @{ E = "  +1.42%  " }
    36 = @{ D = "6.47";      E = "  -0.24%  " }
    37 = @{ D = "2.38";      E = "  +1.42%  " }
    38 = @{ D = "3.60";      E = "  +3.59%  " }
    39 = @{ D = "1.00";      E = "  +0.02%  " }
    40 = @{ D = "18.15";     E = "  -0.95%  " }
    41 = @{ E = "  +4.77%  " }
    42 = @{ D = "1.542.96";  E = "  +0.95%  " }
    43 = @{ D = "100.98";    E = "  +3.37%  " }
    44 = @{ E = "  -0.80%  " }
    45 = @{ D = "0.0920";    E = "  +3.43%  " }
    46 = @{ D = "7.66";      E = "  +8.04%  " }
    47 = @{ E = "  +1.67%  " }
    48 = @{ E = "  -0.27%  " }
    49 = @{ E = "  +3.46%  " }
    50 = @{ D = "2.96";      E = "  +0.55%  " }
    51 = @{ D = "2.291.80";  E = "  +2.42%  " }
}

foreach ($rowNum in $updates.Keys) {
    $cellChanges = $updates[$rowNum]
    foreach ($col in $cellChanges.Keys) {
        $ws.Range("$col$rowNum").Value = $cellChanges[$col]
    }
}
